$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 previously ended with the "footer" block in K9. A new "blog" block
# (ser: 138 - "What is an API") is inserted before it, pushing the footer
# block out to L9.
$footerText = $ws.Cells.Item(9, 11).Value2

$ws.Cells.Item(9, 12).Value = $footerText
$ws.Cells.Item(9, 12).WrapText = $true

$newBlogText = "type: blog" + [char]10 + "width: 2" + [char]10 + "height: 1" + [char]10 + "ser: 138"
$ws.Cells.Item(9, 11).Value = $newBlogText

# Row 8's height shrinks from the max (409.6) down to 255.
$ws.Rows.Item(8).RowHeight = 255

# Selection moves from L8 to K9 (the newly written cell).
$ws.Range("K9").Select()
